# Integrated AG Grid for — portfolio data refresh:
#   * add a "% Change" column (I) to Summary / Joe L / Jonathan R / Michael B / All
#   * refresh existing ticker rows with updated prices / values
#   * append newly-traded tickers (RIOT, MSFT, JBLU) to the relevant sheets
#   * fix a couple of swapped AMD "Sell" transaction rows
#   * append the new RIOT / MSFT / JBLU / NVDA transactions
#
# NOTE: this host mangles COM objects passed through *named* PowerShell
# function parameters (e.g. "Foo -ws $ws"); every helper below is therefore
# called positionally ("Foo $ws ...").

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# Set values on a single row, cell by cell ($null entries are left alone).
function Set-Row {
    param($ws, [int]$row, [object[]]$values)
    for ($i = 0; $i -lt $values.Length; $i++) {
        if ($null -ne $values[$i]) {
            $ws.Cells.Item($row, $i + 1).Value = $values[$i]
        }
    }
}

# Copy the formatting of an existing row onto a brand-new row, then fill in
# its values.
function Add-Row {
    param($ws, [int]$templateRow, [int]$newRow, [int]$lastCol, [object[]]$values)
    $colLetter = [char](64 + $lastCol)
    $ws.Range("A${templateRow}:${colLetter}${templateRow}").Copy() | Out-Null
    $ws.Range("A${newRow}:${colLetter}${newRow}").PasteSpecial($xlPasteFormats) | Out-Null
    Set-Row $ws $newRow $values
}

# Add a new "% Change" header (I1) + data column (I2:I<lastRow>) to a sheet
# whose header/data rows already exist, copying formatting from column H.
function Add-PercentChangeColumn {
    param($ws, [int]$lastRow, [object[]]$values)   # one value per data row, rows 2..lastRow
    $ws.Range("H1").Copy() | Out-Null
    $ws.Range("I1").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Cells.Item(1, 9).Value = "% Change"

    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Range("H$r").Copy() | Out-Null
        $ws.Range("I$r").PasteSpecial($xlPasteFormats) | Out-Null
        $ws.Cells.Item($r, 9).Value = $values[$r - 2]
    }
}

# ===========================================================================
# Summary sheet
# ===========================================================================
$ws = $wb.Worksheets.Item("Summary")

# Updated figures for the existing rows (ticker/quantity/date columns that
# didn't change are left alone).
Set-Row $ws 2  @($null, 812.52, 19,    $null, 15437.88, 8153.1, $null, 7284.78)
Set-Row $ws 3  @($null, 34.66,  $null, $null, 519.9,    $null,  $null, -279.75)
Set-Row $ws 4  @($null, 1285,   $null, $null, 6425,     $null,  $null, 4031)
Set-Row $ws 5  @($null, 230.47, $null, $null, 691.41,   $null,  $null, 513.15)
Set-Row $ws 6  @($null, 26.17,  $null, $null, 130.85,   $null,  $null, 23.25)
Set-Row $ws 7  @($null, 147.22, $null, $null, 1766.64,  $null,  $null, 86.04000000000001)
Set-Row $ws 8  @($null, 151.78, $null, $null, 151.78,   $null,  $null, -3.22)
Set-Row $ws 9  @($null, 332.04, $null, $null, 3984.48,  $null,  $null, -16.2)
Set-Row $ws 10 @($null, 552.58, $null, $null, 2762.9,   $null,  $null, -23.05)
Set-Row $ws 11 @($null, 187.28, $null, $null, 187.28,   $null,  $null, 2.28)

# New ticker rows (RIOT, MSFT, JBLU), formatted like row 11.
Add-Row $ws 11 12 8 @("RIOT", 12.25,  6,  10.55, 73.5,   72.04000000000001, 0.04, 1.46)
Add-Row $ws 11 13 8 @("MSFT", 408.02, 5,  0,     2040.1, 75,                0,    1965.1)
Add-Row $ws 11 14 8 @("JBLU", 5.95,   10, 0,     59.5,   63,                0,    -3.5)

# New "% Change" column across the whole (now 14-row) table.
Add-PercentChangeColumn $ws 14 @(176.1, -35.45, 176.53, 280.11, 22.86, 4.02, -4.1, 0.82, -0.46, -3.69, 6.95, 2573.07, 0)

# ===========================================================================
# Transactions sheet
# ===========================================================================
$ws = $wb.Worksheets.Item("Transactions")

# Rows 11/12 had their Michael B "AMD Sell" entries swapped.
Set-Row $ws 11 @($null, $null, $null, $null, $null, 2, 500, 1000)
Set-Row $ws 12 @($null, $null, $null, $null, $null, 1, 155, 155)

# Newly recorded transactions, formatted like row 16.
Add-Row $ws 16 17 8 @("Jonathan R", 2, 45404, "RIOT", "Buy",  5,  10.51,  52.55)
Add-Row $ws 16 18 8 @("Jonathan R", 2, 45404, "RIOT", "Sell", 1,  10.55,  10.55)
Add-Row $ws 16 19 8 @("Joe L",      2, 45404, "MSFT", "Buy",  5,  15,     75)
Add-Row $ws 16 20 8 @("Jonathan R", 2, 45405, "JBLU", "Buy",  10, 6.3,    63)
Add-Row $ws 16 21 8 @("Jonathan R", 2, 45405, "RIOT", "Buy",  2,  15,     30)
Add-Row $ws 16 22 8 @("Jonathan R", 2, 45405, "NVDA", "Buy",  5,  824.22, 4121.1)

# ===========================================================================
# Joe L sheet
# ===========================================================================
$ws = $wb.Worksheets.Item("Joe L")

Set-Row $ws 2 @($null, 812.52, 19,    $null, 15437.88, 8153.1, $null, 7284.78)
Set-Row $ws 3 @($null, 34.66,  $null, $null, 519.9,    $null,  $null, -279.75)
Set-Row $ws 4 @($null, 1285,   $null, $null, 6425,     $null,  $null, 4031)

# New ticker row (MSFT), formatted like row 4.
Add-Row $ws 4 5 8 @("MSFT", 408.02, 5, 0, 2040.1, 75, 0, 1965.1)

Add-PercentChangeColumn $ws 5 @(176.1, -35.45, 176.53, 2573.07)

# ===========================================================================
# Jonathan R sheet
# ===========================================================================
$ws = $wb.Worksheets.Item("Jonathan R")

Set-Row $ws 2 @($null, 812.52, 19,    $null, 15437.88, 8153.1, $null, 7284.78)
Set-Row $ws 3 @($null, 230.47, $null, $null, 691.41,   $null,  $null, 513.15)
Set-Row $ws 4 @($null, 26.17,  $null, $null, 130.85,   $null,  $null, 23.25)
Set-Row $ws 5 @($null, 147.22, $null, $null, 1766.64,  $null,  $null, 86.04000000000001)

# New ticker rows (RIOT, JBLU), formatted like row 5.
Add-Row $ws 5 6 8 @("RIOT", 12.25, 6,  10.55, 73.5, 72.04000000000001, 0.04, 1.46)
Add-Row $ws 5 7 8 @("JBLU", 5.95,  10, 0,     59.5, 63,                0,    -3.5)

Add-PercentChangeColumn $ws 7 @(176.1, 280.11, 22.86, 4.02, 6.95, 0)

# ===========================================================================
# Michael B sheet
# ===========================================================================
$ws = $wb.Worksheets.Item("Michael B")

Set-Row $ws 2 @($null, 151.78, $null, $null, 151.78,  $null, $null, -3.22)
Set-Row $ws 3 @($null, 332.04, $null, $null, 3984.48, $null, $null, -16.2)
Set-Row $ws 4 @($null, 552.58, $null, $null, 2762.9,  $null, $null, -23.05)
Set-Row $ws 5 @($null, 187.28, $null, $null, 187.28,  $null, $null, 2.28)

Add-PercentChangeColumn $ws 5 @(-4.1, 0.82, -0.46, -3.69)

# ===========================================================================
# All sheet (header only)
# ===========================================================================
$ws = $wb.Worksheets.Item("All")

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(1, 9).Value = "% Change"

Write-Output "done"
